$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Input section ---
# "Change Parameters" (row 5) and "Change Power" (row 6): No -> Yes
$ws.Range("B5").Value = "Yes"
$ws.Range("B6").Value = "Yes"

# --- Output section ---
# "Output .png?" (row 11): Yes -> No
$ws.Range("B11").Value = "No"
# "Output .svg?" (row 12): No -> Yes
$ws.Range("B12").Value = "Yes"
# "Plot Time" (row 17): No -> Yes
$ws.Range("B17").Value = "Yes"
# "Plot Parameter Changes" (row 18): No -> Yes
$ws.Range("B18").Value = "Yes"
# "Plot Power" (row 19): No -> Yes
$ws.Range("B19").Value = "Yes"

# Move the active selection to B4 (was D11)
$ws.Range("B4").Select() | Out-Null
